$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "Adra2b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 150.0354306666667
$ws.Cells.Item(2, 8).Value = 450.106292
$ws.Cells.Item(2, 9).Value = 0.4152507364956075
$ws.Cells.Item(2, 10).Value = 0.4152507364956075
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.002281666666667
$ws.Cells.Item(2, 14).Value = 3.006845
$ws.Cells.Item(2, 15).Value = 0.1474891317421017
$ws.Cells.Item(2, 16).Value = 0.1474891317421018
$ws.Cells.Item(2, 17).Value = 150.3777615076378
$ws.Cells.Item(2, 18).Value = 1353.39985356874
$ws.Cells.Item(2, 19).Value = 0.06124497058100543
$ws.Cells.Item(2, 20).Value = 0.06124497058100543

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "Adra2b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 150.0354306666667
$ws.Cells.Item(3, 8).Value = 450.106292
$ws.Cells.Item(3, 9).Value = 0.4152507364956075
$ws.Cells.Item(3, 10).Value = 0.4152507364956075
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.062332
$ws.Cells.Item(3, 14).Value = 0.186996
$ws.Cells.Item(3, 15).Value = 0.009172364281912123
$ws.Cells.Item(3, 16).Value = 0.009172364281912124
$ws.Cells.Item(3, 17).Value = 9.352008464314666
$ws.Cells.Item(3, 18).Value = 84.168076178832
$ws.Cells.Item(3, 19).Value = 0.003808831023470013
$ws.Cells.Item(3, 20).Value = 0.003808831023470013

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "Adra2b"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 150.0354306666667
$ws.Cells.Item(4, 8).Value = 450.106292
$ws.Cells.Item(4, 9).Value = 0.4152507364956075
$ws.Cells.Item(4, 10).Value = 0.4152507364956075
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.828637333333333
$ws.Cells.Item(4, 14).Value = 8.485911999999999
$ws.Cells.Item(4, 15).Value = 0.4162435353069021
$ws.Cells.Item(4, 16).Value = 0.4162435353069021
$ws.Cells.Item(4, 17).Value = 424.3958205064781
$ws.Cells.Item(4, 18).Value = 3819.562384558304
$ws.Cells.Item(4, 19).Value = 0.1728454345977265
$ws.Cells.Item(4, 20).Value = 0.1728454345977265

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "Adra2b"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 150.0354306666667
$ws.Cells.Item(5, 8).Value = 450.106292
$ws.Cells.Item(5, 9).Value = 0.4152507364956075
$ws.Cells.Item(5, 10).Value = 0.4152507364956075
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.902379666666667
$ws.Cells.Item(5, 14).Value = 8.707139
$ws.Cells.Item(5, 15).Value = 0.4270949686690841
$ws.Cells.Item(5, 16).Value = 0.4270949686690841
$ws.Cells.Item(5, 17).Value = 435.4597832465097
$ws.Cells.Item(5, 18).Value = 3919.138049218588
$ws.Cells.Item(5, 19).Value = 0.1773515002934056
$ws.Cells.Item(5, 20).Value = 0.1773515002934055

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "Adra2b"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 68.382243
$ws.Cells.Item(6, 8).Value = 205.146729
$ws.Cells.Item(6, 9).Value = 0.1892604742946246
$ws.Cells.Item(6, 10).Value = 0.1892604742946246
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.002281666666667
$ws.Cells.Item(6, 14).Value = 3.006845
$ws.Cells.Item(6, 15).Value = 0.1474891317421017
$ws.Cells.Item(6, 16).Value = 0.1474891317421018
$ws.Cells.Item(6, 17).Value = 68.538268484445
$ws.Cells.Item(6, 18).Value = 616.844416360005
$ws.Cells.Item(6, 19).Value = 0.02791386302681255
$ws.Cells.Item(6, 20).Value = 0.02791386302681255

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "Adra2b"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 68.382243
$ws.Cells.Item(7, 8).Value = 205.146729
$ws.Cells.Item(7, 9).Value = 0.1892604742946246
$ws.Cells.Item(7, 10).Value = 0.1892604742946246
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.062332
$ws.Cells.Item(7, 14).Value = 0.186996
$ws.Cells.Item(7, 15).Value = 0.009172364281912123
$ws.Cells.Item(7, 16).Value = 0.009172364281912124
$ws.Cells.Item(7, 17).Value = 4.262401970676
$ws.Cells.Item(7, 18).Value = 38.361617736084
$ws.Cells.Item(7, 19).Value = 0.001735966014397763
$ws.Cells.Item(7, 20).Value = 0.001735966014397762

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "Adra2b"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 68.382243
$ws.Cells.Item(8, 8).Value = 205.146729
$ws.Cells.Item(8, 9).Value = 0.1892604742946246
$ws.Cells.Item(8, 10).Value = 0.1892604742946246
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.828637333333333
$ws.Cells.Item(8, 14).Value = 8.485911999999999
$ws.Cells.Item(8, 15).Value = 0.4162435353069021
$ws.Cells.Item(8, 16).Value = 0.4162435353069021
$ws.Cells.Item(8, 17).Value = 193.428565486872
$ws.Cells.Item(8, 18).Value = 1740.857089381848
$ws.Cells.Item(8, 19).Value = 0.07877844891425563
$ws.Cells.Item(8, 20).Value = 0.0787784489142556

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "Adra2b"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 68.382243
$ws.Cells.Item(9, 8).Value = 205.146729
$ws.Cells.Item(9, 9).Value = 0.1892604742946246
$ws.Cells.Item(9, 10).Value = 0.1892604742946246
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.902379666666667
$ws.Cells.Item(9, 14).Value = 8.707139
$ws.Cells.Item(9, 15).Value = 0.4270949686690841
$ws.Cells.Item(9, 16).Value = 0.4270949686690841
$ws.Cells.Item(9, 17).Value = 198.471231644259
$ws.Cells.Item(9, 18).Value = 1786.241084798331
$ws.Cells.Item(9, 19).Value = 0.08083219633915871
$ws.Cells.Item(9, 20).Value = 0.08083219633915868

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Gnai2"
$ws.Cells.Item(10, 3).Value = "Adra2b"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 104.737245
$ws.Cells.Item(10, 8).Value = 314.211735
$ws.Cells.Item(10, 9).Value = 0.2898796499701289
$ws.Cells.Item(10, 10).Value = 0.2898796499701289
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.002281666666667
$ws.Cells.Item(10, 14).Value = 3.006845
$ws.Cells.Item(10, 15).Value = 0.1474891317421017
$ws.Cells.Item(10, 16).Value = 0.1474891317421018
$ws.Cells.Item(10, 17).Value = 104.976220480675
$ws.Cells.Item(10, 18).Value = 944.7859843260749
$ws.Cells.Item(10, 19).Value = 0.04275409788379868
$ws.Cells.Item(10, 20).Value = 0.04275409788379868

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Gnai2"
$ws.Cells.Item(11, 3).Value = "Adra2b"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 104.737245
$ws.Cells.Item(11, 8).Value = 314.211735
$ws.Cells.Item(11, 9).Value = 0.2898796499701289
$ws.Cells.Item(11, 10).Value = 0.2898796499701289
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.062332
$ws.Cells.Item(11, 14).Value = 0.186996
$ws.Cells.Item(11, 15).Value = 0.009172364281912123
$ws.Cells.Item(11, 16).Value = 0.009172364281912124
$ws.Cells.Item(11, 17).Value = 6.528481955339999
$ws.Cells.Item(11, 18).Value = 58.75633759805999
$ws.Cells.Item(11, 19).Value = 0.002658881747439199
$ws.Cells.Item(11, 20).Value = 0.002658881747439199

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Gnai2"
$ws.Cells.Item(12, 3).Value = "Adra2b"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 104.737245
$ws.Cells.Item(12, 8).Value = 314.211735
$ws.Cells.Item(12, 9).Value = 0.2898796499701289
$ws.Cells.Item(12, 10).Value = 0.2898796499701289
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.828637333333333
$ws.Cells.Item(12, 14).Value = 8.485911999999999
$ws.Cells.Item(12, 15).Value = 0.4162435353069021
$ws.Cells.Item(12, 16).Value = 0.4162435353069021
$ws.Cells.Item(12, 17).Value = 296.26368139748
$ws.Cells.Item(12, 18).Value = 2666.37313257732
$ws.Cells.Item(12, 19).Value = 0.1206605303170938
$ws.Cells.Item(12, 20).Value = 0.1206605303170938

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Gnai2"
$ws.Cells.Item(13, 3).Value = "Adra2b"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 104.737245
$ws.Cells.Item(13, 8).Value = 314.211735
$ws.Cells.Item(13, 9).Value = 0.2898796499701289
$ws.Cells.Item(13, 10).Value = 0.2898796499701289
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.902379666666667
$ws.Cells.Item(13, 14).Value = 8.707139
$ws.Cells.Item(13, 15).Value = 0.4270949686690841
$ws.Cells.Item(13, 16).Value = 0.4270949686690841
$ws.Cells.Item(13, 17).Value = 303.9872502306849
$ws.Cells.Item(13, 18).Value = 2735.885252076165
$ws.Cells.Item(13, 19).Value = 0.1238061400217973
$ws.Cells.Item(13, 20).Value = 0.1238061400217973

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Gnai2"
$ws.Cells.Item(14, 3).Value = "Adra2b"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 38.15794
$ws.Cells.Item(14, 8).Value = 114.47382
$ws.Cells.Item(14, 9).Value = 0.105609139239639
$ws.Cells.Item(14, 10).Value = 0.105609139239639
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.002281666666667
$ws.Cells.Item(14, 14).Value = 3.006845
$ws.Cells.Item(14, 15).Value = 0.1474891317421017
$ws.Cells.Item(14, 16).Value = 0.1474891317421018
$ws.Cells.Item(14, 17).Value = 38.24500369976666
$ws.Cells.Item(14, 18).Value = 344.2050332979
$ws.Cells.Item(14, 19).Value = 0.01557620025048508
$ws.Cells.Item(14, 20).Value = 0.01557620025048508

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Gnai2"
$ws.Cells.Item(15, 3).Value = "Adra2b"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 38.15794
$ws.Cells.Item(15, 8).Value = 114.47382
$ws.Cells.Item(15, 9).Value = 0.105609139239639
$ws.Cells.Item(15, 10).Value = 0.105609139239639
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.062332
$ws.Cells.Item(15, 14).Value = 0.186996
$ws.Cells.Item(15, 15).Value = 0.009172364281912123
$ws.Cells.Item(15, 16).Value = 0.009172364281912124
$ws.Cells.Item(15, 17).Value = 2.37846071608
$ws.Cells.Item(15, 18).Value = 21.40614644472
$ws.Cells.Item(15, 19).Value = 0.0009686854966051489
$ws.Cells.Item(15, 20).Value = 0.0009686854966051488

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Gnai2"
$ws.Cells.Item(16, 3).Value = "Adra2b"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 38.15794
$ws.Cells.Item(16, 8).Value = 114.47382
$ws.Cells.Item(16, 9).Value = 0.105609139239639
$ws.Cells.Item(16, 10).Value = 0.105609139239639
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 2.828637333333333
$ws.Cells.Item(16, 14).Value = 8.485911999999999
$ws.Cells.Item(16, 15).Value = 0.4162435353069021
$ws.Cells.Item(16, 16).Value = 0.4162435353069021
$ws.Cells.Item(16, 17).Value = 107.9349736470933
$ws.Cells.Item(16, 18).Value = 971.41476282384
$ws.Cells.Item(16, 19).Value = 0.04395912147782622
$ws.Cells.Item(16, 20).Value = 0.04395912147782621

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Gnai2"
$ws.Cells.Item(17, 3).Value = "Adra2b"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 38.15794
$ws.Cells.Item(17, 8).Value = 114.47382
$ws.Cells.Item(17, 9).Value = 0.105609139239639
$ws.Cells.Item(17, 10).Value = 0.105609139239639
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 2.902379666666667
$ws.Cells.Item(17, 14).Value = 8.707139
$ws.Cells.Item(17, 15).Value = 0.4270949686690841
$ws.Cells.Item(17, 16).Value = 0.4270949686690841
$ws.Cells.Item(17, 17).Value = 110.7488291778867
$ws.Cells.Item(17, 18).Value = 996.73946260098
$ws.Cells.Item(17, 19).Value = 0.04510513201472256
$ws.Cells.Item(17, 20).Value = 0.04510513201472255
